$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.Value = "'" + $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "69.202.33"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "3.917.61"
$ws.Range("E3").Value = "  +4.70%  "
$ws.Range("E4").Value = "  -0.26%  "
Set-TextValue $ws.Range("D5") "604.31"
$ws.Range("E5").Value = "  +0.63%  "
Set-TextValue $ws.Range("D6") "165.34"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "3.917.57"
$ws.Range("E7").Value = "  +4.76%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue $ws.Range("D9") "0.533"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -1.27%  "
Set-TextValue $ws.Range("D11") "6.40"
$ws.Range("E11").Value = "  -0.06%  "
Set-TextValue $ws.Range("D12") "0.464"
$ws.Range("E12").Value = "  +1.13%  "
Set-TextValue $ws.Range("D13") "37.17"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "4.572.38"
$ws.Range("E15").Value = "  +4.71%  "
$ws.Range("D16").Value = "3.908.68"
$ws.Range("E16").Value = "  +4.07%  "
$ws.Range("D17").Value = "69.278.73"
$ws.Range("E17").Value = "  +0.34%  "
Set-TextValue $ws.Range("D18") "7.52"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("E19").Value = "  -0.32%  "
Set-TextValue $ws.Range("D20") "17.21"
$ws.Range("E20").Value = "  -0.89%  "
Set-TextValue $ws.Range("D21") "11.34"
$ws.Range("E21").Value = "  +2.22%  "
Set-TextValue $ws.Range("D22") "490.43"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("E23").Value = "  -0.06%  "
Set-TextValue $ws.Range("D24") "0.0000167"
$ws.Range("E24").Value = "  +13.45%  "
Set-TextValue $ws.Range("D25") "84.54"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  +0.09%  "
Set-TextValue $ws.Range("D27") "12.21"
$ws.Range("E27").Value = "  -0.88%  "
Set-TextValue $ws.Range("D28") "10.15"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "4.063.77"
$ws.Range("E31").Value = "  +4.38%  "
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D33") "32.60"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D34") "2.39"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("D35").Value = "3.859.40"
$ws.Range("E35").Value = "  +5.06%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  +3.15%  "
Set-TextValue $ws.Range("D38") "0.141"
$ws.Range("E38").Value = "  +2.06%  "
Set-TextValue $ws.Range("D39") "5.94"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D41") "0.321"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D43") "444.60"
$ws.Range("E43").Value = "  +5.00%  "
$ws.Range("E44").Value = "  +0.67%  "
Set-TextValue $ws.Range("D45") "48.49"
$ws.Range("E45").Value = "  -0.06%  "
Set-TextValue $ws.Range("D46") "8.49"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("E47").Value = "  -0.01%  "
Set-TextValue $ws.Range("D48") "27.24"
$ws.Range("E48").Value = "  +16.84%  "
$ws.Range("D49").Value = "2.846.23"
$ws.Range("E49").Value = "  +2.38%  "
Set-TextValue $ws.Range("D50") "141.78"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("E51").Value = "  +2.51%  "
